$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 29999
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""

# Row 23
$ws.Range("H23").Value = 29999
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""

# Row 121
$ws.Range("H121").Value = 2182.5
$ws.Range("I121").Value = 997.5
$ws.Range("J121").Value = 2775
$ws.Range("K121").Value = 2992.5
$ws.Range("L121").Value = 8325
$ws.Range("M121").Value = -1245.5
$ws.Range("N121").Value = -11819

# Row 129
$ws.Range("H129").Value = 922.5
$ws.Range("I129").Value = 577.44446
$ws.Range("K129").Value = 1732.33338
$ws.Range("M129").Value = 3267.66662

# Row 135
$ws.Range("H135").Value = 1105.5769
$ws.Range("I135").Value = 1072.25
$ws.Range("J135").Value = 1216.6666
$ws.Range("K135").Value = 9650.25
$ws.Range("L135").Value = 10949.9994
$ws.Range("M135").Value = -7115.25
$ws.Range("N135").Value = -16019.9994

# Row 137
$ws.Range("H137").Value = 1372.0465
$ws.Range("I137").Value = 1184.3448
$ws.Range("J137").Value = 1760.8572
$ws.Range("K137").Value = 3553.0344
$ws.Range("L137").Value = 5282.571599999999
$ws.Range("M137").Value = -1003.0344
$ws.Range("N137").Value = -10382.5716

# Row 141
$ws.Range("H141").Value = 5885.6597
$ws.Range("I141").Value = 3805.2559
$ws.Range("J141").Value = 28250
$ws.Range("K141").Value = 11415.7677
$ws.Range("L141").Value = 84750
$ws.Range("M141").Value = -6235.7677
$ws.Range("N141").Value = -95110

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1271.7142
$ws.Range("I2").Value = 1150.3334
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1150.3334
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1037.3334
$ws.Range("N2").Value = -2226

# Row 74
$ws.Range("H74").Value = 1238.3846
$ws.Range("I74").Value = 1077.7778
$ws.Range("J74").Value = 1599.75
$ws.Range("K74").Value = 1077.7778
$ws.Range("L74").Value = 1599.75
$ws.Range("M74").Value = -203.7778000000001
$ws.Range("N74").Value = -3347.75

# Row 77
$ws.Range("H77").Value = 1238.3846
$ws.Range("I77").Value = 1077.7778
$ws.Range("J77").Value = 1599.75
$ws.Range("K77").Value = 5388.889
$ws.Range("L77").Value = 7998.75
$ws.Range("M77").Value = -1020.889
$ws.Range("N77").Value = -16734.75

# Row 116
$ws.Range("H116").Value = 1271.7142
$ws.Range("I116").Value = 1150.3334
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1150.3334
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1143.6666
$ws.Range("N116").Value = -6588

# Row 123
$ws.Range("H123").Value = 25222
$ws.Range("J123").Value = 25222
$ws.Range("L123").Value = 25222
$ws.Range("N123").Value = -35022

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1271.7142
$ws.Range("I3").Value = 1150.3334
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1150.3334
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1036.3334
$ws.Range("N3").Value = -2228

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1682.6666
$ws.Range("I16").Value = 1558.2
$ws.Range("J16").Value = 1771.5714
$ws.Range("K16").Value = 1558.2
$ws.Range("L16").Value = 1771.5714
$ws.Range("M16").Value = -1271.2
$ws.Range("N16").Value = -2345.5714

# Row 31
$ws.Range("H31").Value = 2052.7874
$ws.Range("I31").Value = 1462.1945
$ws.Range("K31").Value = 1462.1945
$ws.Range("M31").Value = -1167.1945

# Row 34
$ws.Range("H34").Value = 2052.7874
$ws.Range("I34").Value = 1462.1945
$ws.Range("K34").Value = 1462.1945
$ws.Range("M34").Value = -1260.1945

# Row 35
$ws.Range("H35").Value = 4341.6665
$ws.Range("I35").Value = 4341.6665
$ws.Range("K35").Value = 4341.6665
$ws.Range("M35").Value = -4047.6665

# Row 113
$ws.Range("H113").Value = 1682.6666
$ws.Range("I113").Value = 1558.2
$ws.Range("J113").Value = 1771.5714
$ws.Range("K113").Value = 1558.2
$ws.Range("L113").Value = 1771.5714
$ws.Range("M113").Value = 611.8
$ws.Range("N113").Value = -6111.5714

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 8107.5
$ws.Range("J3").Value = 9960
$ws.Range("L3").Value = 29880
$ws.Range("N3").Value = -30104

# Row 34
$ws.Range("H34").Value = 593.875
$ws.Range("I34").Value = 202
$ws.Range("J34").Value = 620
$ws.Range("K34").Value = 606
$ws.Range("L34").Value = 1860
$ws.Range("M34").Value = -522
$ws.Range("N34").Value = -2028

# Row 129
$ws.Range("H129").Value = 2501306
$ws.Range("I129").Value = 732.1111
$ws.Range("J129").Value = 4547230.5
$ws.Range("K129").Value = 2196.3333
$ws.Range("L129").Value = 13641691.5
$ws.Range("M129").Value = 2803.6667
$ws.Range("N129").Value = -13651691.5

# Row 131
$ws.Range("H131").Value = 16668169
$ws.Range("J131").Value = 16950674
$ws.Range("L131").Value = 50852022
$ws.Range("N131").Value = -50862102

$ws = $wb.Worksheets.Item("GSM")
# Row 110
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180

# Row 123
$ws.Range("H123").Value = 14123.125
$ws.Range("J123").Value = 14123.125
$ws.Range("L123").Value = 14123.125
$ws.Range("N123").Value = -19023.125

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1197.5
$ws.Range("I16").Value = 1197.5
$ws.Range("K16").Value = 1197.5
$ws.Range("M16").Value = -1027.5

# Row 32
$ws.Range("H32").Value = 2494.25
$ws.Range("I32").Value = 2494.25
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2494.25
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""

# Row 61
$ws.Range("H61").Value = 200004
$ws.Range("I61").Value = 200004
$ws.Range("K61").Value = 200004
$ws.Range("M61").Value = -199802

# Row 82
$ws.Range("H82").Value = 1280.3478
$ws.Range("I82").Value = 967.17645
$ws.Range("J82").Value = 2167.6667
$ws.Range("K82").Value = 967.17645
$ws.Range("L82").Value = 2167.6667
$ws.Range("M82").Value = -606.17645
$ws.Range("N82").Value = -2889.6667

# Row 85
$ws.Range("H85").Value = 1280.3478
$ws.Range("I85").Value = 967.17645
$ws.Range("J85").Value = 2167.6667
$ws.Range("K85").Value = 967.17645
$ws.Range("L85").Value = 2167.6667
$ws.Range("M85").Value = 280.82355
$ws.Range("N85").Value = -4663.6667

# Row 100
$ws.Range("H100").Value = 1913.7727
$ws.Range("I100").Value = 1662.6875
$ws.Range("K100").Value = 1662.6875
$ws.Range("M100").Value = -1121.6875

# Row 113
$ws.Range("H113").Value = 200004
$ws.Range("I113").Value = 200004
$ws.Range("K113").Value = 200004
$ws.Range("M113").Value = -197834

$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Range("H28").Value = 15505
$ws.Range("J28").Value = 4006.6667
$ws.Range("L28").Value = 4006.6667
$ws.Range("N28").Value = -4702.6667

# Row 40
$ws.Range("H40").Value = 15000
$ws.Range("J40").Value = 15000
$ws.Range("L40").Value = 15000
$ws.Range("N40").Value = -15298

# Row 64
$ws.Range("H64").Value = 22845.4
$ws.Range("J64").Value = 22845.4
$ws.Range("L64").Value = 22845.4
$ws.Range("N64").Value = -23341.4

# Row 67
$ws.Range("H67").Value = 22845.4
$ws.Range("J67").Value = 22845.4
$ws.Range("L67").Value = 22845.4
$ws.Range("N67").Value = -24561.4

# Row 123
$ws.Range("H123").Value = 21603.842
$ws.Range("J123").Value = 21603.842
$ws.Range("L123").Value = 21603.842
$ws.Range("N123").Value = -31403.842
